$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '97.437.59'
Set-TextValue 'E2' '  -1.36%  '
Set-TextValue 'D3' '3.331.04'
Set-TextValue 'E3' '  -2.58%  '
Set-TextValue 'E4' '  +0.11%  '
Set-TextValue 'D5' '247.79'
Set-TextValue 'D6' '653.99'
Set-TextValue 'E6' '  -2.51%  '
Set-TextValue 'E7' '  -11.47%  '
Set-TextValue 'D8' '0.418'
Set-TextValue 'E8' '  -11.19%  '
Set-TextValue 'D9' '1.00'
Set-TextValue 'E9' '  +0.06%  '
Set-TextValue 'D10' '1.00'
Set-TextValue 'E10' '  -8.39%  '
Set-TextValue 'D11' '3.328.91'
Set-TextValue 'E11' '  -2.53%  '
Set-TextValue 'E12' '  -6.67%  '
Set-TextValue 'D13' '40.28'
Set-TextValue 'E13' '  -6.93%  '
Set-TextValue 'D14' '97.333.64'
Set-TextValue 'E14' '  -1.60%  '
Set-TextValue 'D15' '6.04'
Set-TextValue 'E15' '  -1.34%  '
Set-TextValue 'D16' '0.0000252'
Set-TextValue 'E16' '  -9.04%  '
Set-TextValue 'D17' '3.955.92'
Set-TextValue 'E17' '  -2.64%  '
Set-TextValue 'D18' '8.51'
Set-TextValue 'E18' '  +4.60%  '
Set-TextValue 'D19' '3.332.39'
Set-TextValue 'E19' '  -2.57%  '
Set-TextValue 'D20' '0.533'
Set-TextValue 'E20' '  +21.51%  '
Set-TextValue 'D21' '16.74'
Set-TextValue 'E21' '  -3.52%  '
Set-TextValue 'D22' '10.59'
Set-TextValue 'E22' '  -1.10%  '
Set-TextValue 'D23' '496.54'
Set-TextValue 'E23' '  -7.79%  '
Set-TextValue 'D24' '3.28'
Set-TextValue 'E24' '  -8.60%  '
Set-TextValue 'D25' '0.0000198'
Set-TextValue 'E25' '  -9.75%  '
Set-TextValue 'D26' '6.45'
Set-TextValue 'E26' '  +0.58%  '
Set-TextValue 'D27' '93.27'
Set-TextValue 'E27' '  -9.66%  '
Set-TextValue 'D28' '11.98'
Set-TextValue 'E28' '  -7.04%  '
Set-TextValue 'D29' '3.512.63'
Set-TextValue 'E29' '  -2.40%  '
Set-TextValue 'D30' '0.145'
Set-TextValue 'E30' '  -4.07%  '
Set-TextValue 'E31' '  -0.20%  '
Set-TextValue 'D32' '10.87'
Set-TextValue 'E32' '  -7.03%  '
Set-TextValue 'E33' '  -5.18%  '
Set-TextValue 'D34' '2.47'
Set-TextValue 'E34' '  +11.38%  '
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  -0.91%  '
Set-TextValue 'D36' '0.547'
Set-TextValue 'E36' '  -3.36%  '
Set-TextValue 'D37' '28.27'
Set-TextValue 'E37' '  -6.84%  '
Set-TextValue 'D38' '7.54'
Set-TextValue 'E38' '  -5.11%  '
Set-TextValue 'D39' '1.43'
Set-TextValue 'E39' '  +0.44%  '
Set-TextValue 'E40' '  -0.05%  '
Set-TextValue 'E41' '  -7.55%  '
Set-TextValue 'D42' '503.32'
Set-TextValue 'E42' '  -5.70%  '
Set-TextValue 'D43' '24.59'
Set-TextValue 'E43' '  -0.79%  '
Set-TextValue 'E44' '  -2.81%  '
Set-TextValue 'D45' '0.833'
Set-TextValue 'E45' '  -3.14%  '
Set-TextValue 'D46' '8.59'
Set-TextValue 'E46' '  +4.60%  '
Set-TextValue 'D47' '0.0409'
Set-TextValue 'E47' '  -7.04%  '
Set-TextValue 'D48' '5.45'
Set-TextValue 'E48' '  +2.15%  '
Set-TextValue 'D49' '1.63'
Set-TextValue 'E49' '  +2.46%  '
Set-TextValue 'D50' '53.37'
Set-TextValue 'E50' '  +5.18%  '
Set-TextValue 'D51' '3.13'
Set-TextValue 'E51' '  -11.65%  '
